$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell A3 value from 0 to 1
$ws.Range("A3").Value = 1

# Scroll the view back so A1 is the top-left visible cell again (was D1)
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1

# Move the active selection to A3 (was J11)
$ws.Range("A3").Select()
